$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 has the exact same style pattern (A=date style, B=plain centered, C=wrapped)
# that the new row 21 needs, so copy its formatting first, then overwrite values.
$ws.Range("A19:C19").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A21").Value = 45246
$ws.Range("B21").Value = "~2 hrs"
$ws.Range("C21").Value = "added the syntax highlighter, fixed some issues with negative numbers, and finished the print result section"
$ws.Rows.Item(21).RowHeight = 45

# Row 22 only has a date in column A, formatted like the other date cells (A19 style).
$ws.Range("A19").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A22").Value = 45247

# Match the new selection recorded in the workbook (active cell C22).
$ws.Range("C22").Select()
